# Auto-generated script to apply scheduled-runner market data updates
# to the Ultima_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3193.6875
$ws.Range("I51").Value = 2400
$ws.Range("J51").Value = 3376.8462
$ws.Range("K51").Value = 2400
$ws.Range("L51").Value = 3376.8462
$ws.Range("M51").Value = -1916
$ws.Range("N51").Value = -4344.8462
$ws.Range("H64").Value = 3391.6667
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 3320
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 3320
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -3816
$ws.Range("H67").Value = 3391.6667
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 3320
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 3320
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -5036
$ws.Range("H98").Value = 1581.5264
$ws.Range("I98").Value = 1179.1666
$ws.Range("J98").Value = 2271.2856
$ws.Range("K98").Value = 1179.1666
$ws.Range("L98").Value = 2271.2856
$ws.Range("M98").Value = 318.8334
$ws.Range("N98").Value = -5267.2856
$ws.Range("H122").Value = 1581.5264
$ws.Range("I122").Value = 1179.1666
$ws.Range("J122").Value = 2271.2856
$ws.Range("K122").Value = 3537.4998
$ws.Range("L122").Value = 6813.8568
$ws.Range("M122").Value = -1087.4998
$ws.Range("N122").Value = -11713.8568
$ws.Range("H137").Value = 12501894
$ws.Range("J137").Value = 20002252
$ws.Range("L137").Value = 60006756
$ws.Range("N137").Value = -60011856
$ws.Range("H138").Value = 2834.2942
$ws.Range("J138").Value = 3866.6667
$ws.Range("L138").Value = 11600.0001
$ws.Range("N138").Value = -21880.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10016.667
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 50
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = -15230
$ws.Range("H7").Value = 50712
$ws.Range("J7").Value = 50712
$ws.Range("L7").Value = 50712
$ws.Range("N7").Value = -50940
$ws.Range("H8").Value = 80000
$ws.Range("J8").Value = 80000
$ws.Range("L8").Value = 80000
$ws.Range("N8").Value = -80288
$ws.Range("H61").Value = 16130939
$ws.Range("I61").Value = 17859056
$ws.Range("J61").Value = 1838
$ws.Range("K61").Value = 17859056
$ws.Range("L61").Value = 1838
$ws.Range("M61").Value = -17858844
$ws.Range("N61").Value = -2262
$ws.Range("H74").Value = 17860862
$ws.Range("I74").Value = 33335736
$ws.Range("J74").Value = 5238.231
$ws.Range("K74").Value = 33335736
$ws.Range("L74").Value = 5238.231
$ws.Range("M74").Value = -33334862
$ws.Range("N74").Value = -6986.231
$ws.Range("H77").Value = 17860862
$ws.Range("I77").Value = 33335736
$ws.Range("J77").Value = 5238.231
$ws.Range("K77").Value = 166678680
$ws.Range("L77").Value = 26191.155
$ws.Range("M77").Value = -166674312
$ws.Range("N77").Value = -34927.155
$ws.Range("H97").Value = 7389.125
$ws.Range("I97").Value = 9762
$ws.Range("J97").Value = 2168.8
$ws.Range("K97").Value = 9762
$ws.Range("L97").Value = 2168.8
$ws.Range("M97").Value = -9266
$ws.Range("N97").Value = -3160.8
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H132").Value = 4547009.5
$ws.Range("I132").Value = 5556829.5
$ws.Range("K132").Value = 16670488.5
$ws.Range("M132").Value = -16667958.5
$ws.Range("H136").Value = 16130939
$ws.Range("I136").Value = 17859056
$ws.Range("J136").Value = 1838
$ws.Range("K136").Value = 53577168
$ws.Range("L136").Value = 5514
$ws.Range("M136").Value = -53574618
$ws.Range("N136").Value = -10614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1373.5416
$ws.Range("I94").Value = 1014.3333
$ws.Range("J94").Value = 1972.2222
$ws.Range("K94").Value = 1014.3333
$ws.Range("L94").Value = 1972.2222
$ws.Range("M94").Value = -563.3333
$ws.Range("N94").Value = -2874.2222
$ws.Range("H134").Value = 3318.1943
$ws.Range("I134").Value = 1959.1
$ws.Range("J134").Value = 5017.0625
$ws.Range("K134").Value = 5877.299999999999
$ws.Range("L134").Value = 15051.1875
$ws.Range("M134").Value = -3342.299999999999
$ws.Range("N134").Value = -20121.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -387
$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5340
$ws.Range("H31").Value = 11910915
$ws.Range("I31").Value = 11120.923
$ws.Range("J31").Value = 22224070
$ws.Range("K31").Value = 11120.923
$ws.Range("L31").Value = 22224070
$ws.Range("M31").Value = -10825.923
$ws.Range("N31").Value = -22224660
$ws.Range("H34").Value = 11910915
$ws.Range("I34").Value = 11120.923
$ws.Range("J34").Value = 22224070
$ws.Range("K34").Value = 11120.923
$ws.Range("L34").Value = 22224070
$ws.Range("M34").Value = -10918.923
$ws.Range("N34").Value = -22224474
$ws.Range("H62").Value = 2328.5715
$ws.Range("I62").Value = 2328.5715
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2328.5715
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1704.5715
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2328.5715
$ws.Range("I65").Value = 2328.5715
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11642.8575
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -8522.8575
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 2151.5
$ws.Range("I99").Value = 1199.5
$ws.Range("J99").Value = 3103.5
$ws.Range("K99").Value = 1199.5
$ws.Range("L99").Value = 3103.5
$ws.Range("M99").Value = 298.5
$ws.Range("N99").Value = -6099.5
$ws.Range("H126").Value = 2151.5
$ws.Range("I126").Value = 1199.5
$ws.Range("J126").Value = 3103.5
$ws.Range("K126").Value = 3598.5
$ws.Range("L126").Value = 9310.5
$ws.Range("M126").Value = -1128.5
$ws.Range("N126").Value = -14250.5
$ws.Range("H132").Value = 6074
$ws.Range("I132").Value = 6142.1904
$ws.Range("J132").Value = 5835.3335
$ws.Range("K132").Value = 18426.5712
$ws.Range("L132").Value = 17506.0005
$ws.Range("M132").Value = -15896.5712
$ws.Range("N132").Value = -22566.0005
$ws.Range("H134").Value = 2552.1667
$ws.Range("I134").Value = 2563.4614
$ws.Range("J134").Value = 2522.8
$ws.Range("K134").Value = 7690.3842
$ws.Range("L134").Value = 7568.400000000001
$ws.Range("M134").Value = -5155.3842
$ws.Range("N134").Value = -12638.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 6008.727
$ws.Range("J96").Value = 6677.3335
$ws.Range("L96").Value = 20032.0005
$ws.Range("N96").Value = -24150.0005
$ws.Range("H113").Value = 1349.091
$ws.Range("I113").Value = 480
$ws.Range("J113").Value = 3666.6667
$ws.Range("K113").Value = 1440
$ws.Range("L113").Value = 11000.0001
$ws.Range("M113").Value = 730
$ws.Range("N113").Value = -15340.0001
$ws.Range("H120").Value = 12792.556
$ws.Range("I120").Value = 4766.6665
$ws.Range("J120").Value = 16805.5
$ws.Range("K120").Value = 14299.9995
$ws.Range("L120").Value = 50416.5
$ws.Range("M120").Value = -9461.999500000002
$ws.Range("N120").Value = -60092.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1445.1364
$ws.Range("I97").Value = 1255.9333
$ws.Range("J97").Value = 1850.5714
$ws.Range("K97").Value = 1255.9333
$ws.Range("L97").Value = 1850.5714
$ws.Range("M97").Value = -759.9332999999999
$ws.Range("N97").Value = -2842.5714
$ws.Range("H122").Value = 6062170.5
$ws.Range("I122").Value = 7408652
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 22225956
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -22223506
$ws.Range("N122").Value = -13912
$ws.Range("H132").Value = 5997.4165
$ws.Range("I132").Value = 5366.5
$ws.Range("K132").Value = 16099.5
$ws.Range("M132").Value = -13569.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11913151
$ws.Range("I132").Value = 5851.625
$ws.Range("J132").Value = 27789550
$ws.Range("K132").Value = 17554.875
$ws.Range("L132").Value = 83368650
$ws.Range("M132").Value = -15024.875
$ws.Range("N132").Value = -83373710
$ws.Range("H136").Value = 41679576
$ws.Range("I136").Value = 55557430
$ws.Range("K136").Value = 166672290
$ws.Range("M136").Value = -166669740
$ws.Range("H139").Value = 44941.168
$ws.Range("J139").Value = 45331.273
$ws.Range("L139").Value = 45331.273
$ws.Range("N139").Value = -55611.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 37320
$ws.Range("J75").Value = 37320
$ws.Range("L75").Value = 37320
$ws.Range("N75").Value = -39192
$ws.Range("H78").Value = 37320
$ws.Range("J78").Value = 37320
$ws.Range("L78").Value = 111960
$ws.Range("N78").Value = -121320
$ws.Range("H132").Value = 1269.9016
$ws.Range("I132").Value = 995.4792
$ws.Range("K132").Value = 2986.4376
$ws.Range("M132").Value = -456.4376000000002
$ws.Range("H136").Value = 1300.25
$ws.Range("I136").Value = 1228.8572
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 3686.5716
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -1136.5716
$ws.Range("N136").Value = -10500

